$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Handoff / HO Xliff generation timestamps ---
$wsOverview.Range("G2").Value = "2016-09-05 14:49:43"
$wsZhCn.Range("H2").Value = "2016-09-05 14:49:39"
$wsDeDe.Range("H2").Value = "2016-09-05 14:49:43"

# --- Widen the status-related columns to fit "Ready for handoff" ---
# (input value chosen so the persisted OOXML column width rounds to
#  17.166666666666668, the closest value to the target 17.2159881591797
#  reachable through the ColumnWidth property's internal pixel rounding)
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333336
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333336
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333336
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333336
